# Insert a new price record as row 52 in the "Achicoria" sheet, pushing the
# existing rows 52:124 down to 53:125 (dimension grows from A1:R124 to A1:R125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52:124 down by one row, carrying their formatting with them.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly record.
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = 45117
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = 100112010
$ws.Range("G52").Value = "Achicoria"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 150
$ws.Range("K52").Value = 8000
$ws.Range("L52").Value = 8000
$ws.Range("M52").Value = 8000
$ws.Range("N52").Value = "$/caja 18 unidades"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 444
$ws.Range("Q52").Value = 18
$ws.Range("R52").Value = "Hortaliza"
